$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.984.26'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '3.414.71'
$ws.Range('E3').Value = '  -2.81%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '406.31'
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '135.06'
$ws.Range('E6').Value = '  +4.41%  '
$ws.Range('E7').Value = '  -1.31%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.686'
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('E10').Value = '  -6.71%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '42.78'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '8.43'
$ws.Range('E13').Value = '  -3.67%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '19.92'
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').Value = '3.410.55'
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').Value = '61.959.30'
$ws.Range('E16').Value = '  -1.75%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.02'
$ws.Range('E17').Value = '  -3.35%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.02'
$ws.Range('E18').Value = '  -1.72%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000130'
$ws.Range('E19').Value = '  -5.51%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.19'
$ws.Range('E20').Value = '  -5.28%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '84.75'
$ws.Range('E21').Value = '  +3.08%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '312.91'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.90'
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('E25').Value = '  +8.82%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '29.62'
$ws.Range('E26').Value = '  -3.29%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.15'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('E28').Value = '  +4.31%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.59'
$ws.Range('E29').Value = '  -2.77%  '
$ws.Range('E30').Value = '  -4.59%  '
$ws.Range('E31').Value = '  -2.96%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '42.40'
$ws.Range('E32').Value = '  -2.82%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.32'
$ws.Range('E34').Value = '  -6.63%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0482'
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '51.74'
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.41'
$ws.Range('E38').Value = '  -4.87%  '
$ws.Range('E39').Value = '  -2.93%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.306'
$ws.Range('E40').Value = '  +6.01%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '137.44'
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.98'
$ws.Range('E42').Value = '  -1.78%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.125'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.04'
$ws.Range('E44').Value = '  +1.36%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '16.70'
$ws.Range('E45').Value = '  -6.26%  '
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '21.29'
$ws.Range('E47').Value = '  -5.06%  '
$ws.Range('D48').Value = '2.117.20'
$ws.Range('E48').Value = '  -4.80%  '
$ws.Range('E49').Value = '  -3.25%  '
$ws.Range('E50').Value = '  +0.43%  '
$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0349'
$ws.Range('E51').Value = '  +2.80%  '
